$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from O1 into P1:Q1 so the new header cells get the bold/border/center style
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 1 new header values
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: swap I/K and M/O values, and fill in new P/Q columns with 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}

Write-Output "done"
